$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting -------------------------------------------------------
# Rows 77-80 (new WVU Stadium monitoring sites) reuse the style family
# already used by row 40 (s=39/47/48/49 across columns A,C,D,G,H,I,J,M,N,
# O,P,Q,R,S,T) except column G, which instead matches the s=42 family used
# by rows 33/36/37. Copy formats cell-by-cell so only the columns that
# actually carry data in the new rows end up with an explicit style
# (matching the target layout, which leaves B,E,F,K,L,U completely empty).
$formatCols = @("A","C","D","H","I","J","M","N","O","P","Q","R","S","T")
foreach ($targetRow in 77..80) {
  foreach ($col in $formatCols) {
    $ws.Range("$col" + "40").Copy() | Out-Null
    $ws.Range("$col" + "$targetRow").PasteSpecial(-4122) | Out-Null
  }
  $ws.Range("G33").Copy() | Out-Null
  $ws.Range("G$targetRow").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- Values -------------------------------------------------------------
# Row 77 - WVU Stadium SW
$ws.Range("A77").Value = "StadiumSW-01"
$ws.Range("C77").Value = "active"
$ws.Range("D77").Value = "WVU Stadium SW"
$ws.Range("G77").Value = "not institution specific"
$ws.Range("H77").Value = "upstream"
$ws.Range("I77").Value = "Sewer Network"
$ws.Range("J77").Value = "Campus"
$ws.Range("M77").Value = "StarCityWWTP-01"
$ws.Range("N77").Value = "Monongalia"
$ws.Range("Q77").Value = 105612
$ws.Range("R77").Value = "12-hr time-weighted composite"
$ws.Range("S77").Value = "raw wastewater"
$ws.Range("T77").Value = 26505

# Rows 78-80: column A then column D were filled in across all three rows
# before the remaining columns (matches the shared-string append order in
# the saved workbook).
$ws.Range("A78").Value = "StadiumNW-01"
$ws.Range("A79").Value = "StadiumNE-01"
$ws.Range("A80").Value = "StadiumSE-01"
$ws.Range("D78").Value = "WVU Stadium NW"
$ws.Range("D79").Value = "WVU Stadium NE"
$ws.Range("D80").Value = "WVU Stadium SE"

# Row 78 - WVU Stadium NW
$ws.Range("C78").Value = "active"
$ws.Range("G78").Value = "not institution specific"
$ws.Range("H78").Value = "upstream"
$ws.Range("I78").Value = "Sewer Network"
$ws.Range("J78").Value = "Campus"
$ws.Range("M78").Value = "StarCityWWTP-01"
$ws.Range("N78").Value = "Monongalia"
$ws.Range("Q78").Value = 105612
$ws.Range("R78").Value = "12-hr time-weighted composite"
$ws.Range("S78").Value = "raw wastewater"
$ws.Range("T78").Value = 26505

# Row 79 - WVU Stadium NE
$ws.Range("C79").Value = "active"
$ws.Range("G79").Value = "not institution specific"
$ws.Range("H79").Value = "upstream"
$ws.Range("I79").Value = "Sewer Network"
$ws.Range("J79").Value = "Campus"
$ws.Range("M79").Value = "StarCityWWTP-01"
$ws.Range("N79").Value = "Monongalia"
$ws.Range("Q79").Value = 105612
$ws.Range("R79").Value = "12-hr time-weighted composite"
$ws.Range("S79").Value = "raw wastewater"
$ws.Range("T79").Value = 26505

# Row 80 - WVU Stadium SE
$ws.Range("C80").Value = "active"
$ws.Range("G80").Value = "not institution specific"
$ws.Range("H80").Value = "upstream"
$ws.Range("I80").Value = "Sewer Network"
$ws.Range("J80").Value = "Campus"
$ws.Range("M80").Value = "StarCityWWTP-01"
$ws.Range("N80").Value = "Monongalia"
$ws.Range("Q80").Value = 105612
$ws.Range("R80").Value = "12-hr time-weighted composite"
$ws.Range("S80").Value = "raw wastewater"
$ws.Range("T80").Value = 26505

# --- View state -----------------------------------------------------
# Update the frozen-pane scroll position and final selection to match
# the author's saved view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 51
$win.ScrollColumn = 2
$ws.Range("G80").Select()
